# "Atualiza instrucao de trabalho"
#
# - Flip Status (col J) from "Pendente" to "Resolvido" for 4 incidents on
#   the ITI sheet (rows 9, 17, 19, 22) that just got resolved.
# - The sheet's AutoFilter used to hide everything except the "Pendente"
#   rows for a given month; that filter is cleared so every row on ITI is
#   visible again (keeps the AutoFilter range + existing sort, just drops
#   the filter criteria and un-hides the rows it was hiding).
# - Selection/cursor position bookkeeping on both sheets.

$wb  = $excel.ActiveWorkbook
$iti = $wb.Worksheets.Item("ITI")
$spn = $wb.Worksheets.Item("SPN")

# --- Status column updates: Pendente -> Resolvido ------------------------
$iti.Range("J9").Value  = "Resolvido"
$iti.Range("J17").Value = "Resolvido"
$iti.Range("J19").Value = "Resolvido"
$iti.Range("J22").Value = "Resolvido"

# --- Clear the AutoFilter criteria and unhide the filtered-out rows ------
# (keeps the AutoFilter range itself, and the sortState/sortCondition)
$iti.ShowAllData()

# --- Window size/position bookkeeping (best effort) -----------------------
try {
  $win = $excel.Windows.Item(1)
  $win.Left   = -120
  $win.Top    = -120
  $win.Width  = 29040
  $win.Height = 15720
} catch {}

# --- Selection bookkeeping -------------------------------------------------
# Set SPN's selection first, then re-activate ITI and select there last so
# ITI (the tab that was already active) stays the active/visible tab.
$spn.Range("J3").Select() | Out-Null
$iti.Activate() | Out-Null
$iti.Range("F20").Select() | Out-Null
